$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> (B, C, D, E, F)
$data = @{
  2  = @(1,  0.03958333333333333, 0.06180555555555556, "E4",  1)
  3  = @(2,  0.01388888888888889, 0.03958333333333333, "D48", 1)
  4  = @(3,  0.03402777777777777, 0.05625,              "C25", 3)
  5  = @(4,  0.03402777777777777, 0.05694444444444444, "A13", 1)
  6  = @(5,  0.02847222222222222, 0.05416666666666667, "E2",  1)
  7  = @(6,  0,                   0.02569444444444444, "A1",  1)
  8  = @(7,  0.03958333333333333, 0.06597222222222222, "A18", 2)
  9  = @(8,  0.003472222222222222,0.02986111111111111, "A11", 3)
  10 = @(9,  0.03819444444444445, 0.06180555555555556, "F40", 2)
  11 = @(10, 0.03888888888888889, 0.06597222222222222, "F59", 3)
  12 = @(11, 0.04097222222222222, 0.06666666666666667, "C3",  2)
  13 = @(12, 0.009027777777777777,0.03402777777777777, "B8",  1)
  14 = @(13, 0.0125,              0.03888888888888889, "C26", 2)
  15 = @(14, 0.01805555555555555, 0.04097222222222222, "A2",  1)
  16 = @(15, 0.03333333333333333, 0.05694444444444444, "D46", 1)
  17 = @(16, 0.03055555555555555, 0.05416666666666667, "D42", 3)
  18 = @(17, 0.03263888888888889, 0.05902777777777778, "A14", 1)
  19 = @(18, 0.003472222222222222,0.025,                "C15", 2)
  20 = @(19, 0.01805555555555555, 0.04166666666666666, "D41", 3)
  21 = @(20, 0.02777777777777778, 0.04930555555555555, "D32", 3)
}

foreach ($row in $data.Keys) {
  $vals = $data[$row]
  $ws.Cells.Item($row, 2).Value = $vals[0]
  $ws.Cells.Item($row, 3).Value = $vals[1]
  $ws.Cells.Item($row, 4).Value = $vals[2]
  $ws.Cells.Item($row, 5).Value = $vals[3]
  $ws.Cells.Item($row, 6).Value = $vals[4]
}
